$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-13 from
# serial date 45185 (2023-09-16) to 45204 (2023-10-05).
$newDate = Get-Date -Year 2023 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Range("C2:C13").Value = $newDate
